# Applies commit '#5: insurance, claim, debt, investment done':
# adds property_category/category/date/legislator_name/legislator_id/source_file/index
# metadata columns (matching the pattern already used on the other sheets) to the
# insurance (sheet 6), claim (sheet 7) and investment (sheet 8) tables, and turns their
# header rows into proper field-name labels instead of duplicated first-row data.
$wb = $excel.ActiveWorkbook

# ---- sheet 6 ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").Value = "2012-04-30"
$ws.Range("H2").Value = "鄭麗君"
$ws.Range("I2").Value = 1764
$ws.Range("J2").Value = "tmp81521"
$ws.Range("K2").Value = 114
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").Value = "2012-04-30"
$ws.Range("H3").Value = "鄭麗君"
$ws.Range("I3").Value = 1764
$ws.Range("J3").Value = "tmp81521"
$ws.Range("K3").Value = 115
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
$ws.Range("G4").Value = "2012-04-30"
$ws.Range("H4").Value = "鄭麗君"
$ws.Range("I4").Value = 1764
$ws.Range("J4").Value = "tmp81521"
$ws.Range("K4").Value = 116
$ws.Range("E1:K1").Font.Bold = $true
$ws.Range("E1:K1").Borders.LineStyle = 1
$ws.Range("E1:K1").HorizontalAlignment = -4108
$ws.Range("E1:K1").VerticalAlignment = -4160

# ---- sheet 7 ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "debtor"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
$ws.Range("H2").Value = "claim"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-30"
$ws.Range("K2").Value = "鄭麗君"
$ws.Range("L2").Value = 1764
$ws.Range("M2").Value = "tmp81521"
$ws.Range("N2").Value = 121
$ws.Range("H3").Value = "claim"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-04-30"
$ws.Range("K3").Value = "鄭麗君"
$ws.Range("L3").Value = 1764
$ws.Range("M3").Value = "tmp81521"
$ws.Range("N3").Value = 122
$ws.Range("H4").Value = "claim"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2012-04-30"
$ws.Range("K4").Value = "鄭麗君"
$ws.Range("L4").Value = 1764
$ws.Range("M4").Value = "tmp81521"
$ws.Range("N4").Value = 124
$ws.Range("H5").Value = "claim"
$ws.Range("I5").Value = "normal"
$ws.Range("J5").Value = "2012-04-30"
$ws.Range("K5").Value = "鄭麗君"
$ws.Range("L5").Value = 1764
$ws.Range("M5").Value = "tmp81521"
$ws.Range("N5").Value = 125
$ws.Range("H6").Value = "claim"
$ws.Range("I6").Value = "normal"
$ws.Range("J6").Value = "2012-04-30"
$ws.Range("K6").Value = "鄭麗君"
$ws.Range("L6").Value = 1764
$ws.Range("M6").Value = "tmp81521"
$ws.Range("N6").Value = 126
$ws.Range("H7").Value = "claim"
$ws.Range("I7").Value = "normal"
$ws.Range("J7").Value = "2012-04-30"
$ws.Range("K7").Value = "鄭麗君"
$ws.Range("L7").Value = 1764
$ws.Range("M7").Value = "tmp81521"
$ws.Range("N7").Value = 127
$ws.Range("H8").Value = "claim"
$ws.Range("I8").Value = "normal"
$ws.Range("J8").Value = "2012-04-30"
$ws.Range("K8").Value = "鄭麗君"
$ws.Range("L8").Value = 1764
$ws.Range("M8").Value = "tmp81521"
$ws.Range("N8").Value = 128
$ws.Range("H1:N1").Font.Bold = $true
$ws.Range("H1:N1").Borders.LineStyle = 1
$ws.Range("H1:N1").HorizontalAlignment = -4108
$ws.Range("H1:N1").VerticalAlignment = -4160

# ---- sheet 8 ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("B1").Value = "owner"
$ws.Range("C1").Value = "company"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
$ws.Range("H2").Value = "investment"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-30"
$ws.Range("K2").Value = "鄭麗君"
$ws.Range("L2").Value = 1764
$ws.Range("M2").Value = "tmp81521"
$ws.Range("N2").Value = 138
$ws.Range("H3").Value = "investment"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-04-30"
$ws.Range("K3").Value = "鄭麗君"
$ws.Range("L3").Value = 1764
$ws.Range("M3").Value = "tmp81521"
$ws.Range("N3").Value = 139
$ws.Range("H4").Value = "investment"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2012-04-30"
$ws.Range("K4").Value = "鄭麗君"
$ws.Range("L4").Value = 1764
$ws.Range("M4").Value = "tmp81521"
$ws.Range("N4").Value = 140
$ws.Range("H5").Value = "investment"
$ws.Range("I5").Value = "normal"
$ws.Range("J5").Value = "2012-04-30"
$ws.Range("K5").Value = "鄭麗君"
$ws.Range("L5").Value = 1764
$ws.Range("M5").Value = "tmp81521"
$ws.Range("N5").Value = 141
$ws.Range("H6").Value = "investment"
$ws.Range("I6").Value = "normal"
$ws.Range("J6").Value = "2012-04-30"
$ws.Range("K6").Value = "鄭麗君"
$ws.Range("L6").Value = 1764
$ws.Range("M6").Value = "tmp81521"
$ws.Range("N6").Value = 142
$ws.Range("H7").Value = "investment"
$ws.Range("I7").Value = "normal"
$ws.Range("J7").Value = "2012-04-30"
$ws.Range("K7").Value = "鄭麗君"
$ws.Range("L7").Value = 1764
$ws.Range("M7").Value = "tmp81521"
$ws.Range("N7").Value = 143
$ws.Range("H1:N1").Font.Bold = $true
$ws.Range("H1:N1").Borders.LineStyle = 1
$ws.Range("H1:N1").HorizontalAlignment = -4108
$ws.Range("H1:N1").VerticalAlignment = -4160

